$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# This edit simulates a failed handoff transform for the file that used to
# be "a6e1d776-ca2c-4c1e-b047-befc8dde8234.md": it is renamed to
# "66f80cc9-80c1-471b-b042-8610d0a171d2.md", its status flips from
# "Ready for handoff" to "Handoff transform failed", and on the per-locale
# sheets the (now nonexistent) handoff artifact/date info is cleared out,
# with the handoff reason moving from "Include" to "Ignored".
# -----------------------------------------------------------------------

$oldName = "a6e1d776-ca2c-4c1e-b047-befc8dde8234.md"
$newName = "66f80cc9-80c1-471b-b042-8610d0a171d2.md"
$oldStatus = "Ready for handoff"
$newStatus = "Handoff transform failed"
$epoch = "0001-01-01 00:00:00"
$newReason = "Ignored"

function Update-Md-Hyperlink($ws) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.TextToDisplay -eq $oldName) {
            $h.TextToDisplay = $newName
        }
    }
}

# ---------------- Overview sheet ----------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newName
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
Update-Md-Hyperlink $wsOverview

# ---------------- zh-cn sheet ----------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newName
$wsZh.Range("B2").Value = $newStatus

# The handoff for this locale failed, so the "Latest Handoff File" cell no
# longer has any content or hyperlink.
foreach ($h in $wsZh.Hyperlinks) {
    if ($h.TextToDisplay -eq "a6e1d776-ca2c-4c1e-b047-befc8dde8234.4fd030535bba4f8016cf207827d76bef0020e7fb.zh-cn.xlf") {
        $h.Delete()
    }
}
$wsZh.Range("C2").Clear()

$wsZh.Range("D2").Value = $epoch
$wsZh.Range("G2").Value = $epoch
$wsZh.Range("H2").Value = $newReason

$wsZh.Range("D3").Value = $epoch
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = $newReason

Update-Md-Hyperlink $wsZh

# ---------------- de-de sheet ----------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newName
$wsDe.Range("B2").Value = $newStatus

foreach ($h in $wsDe.Hyperlinks) {
    if ($h.TextToDisplay -eq "a6e1d776-ca2c-4c1e-b047-befc8dde8234.4fd030535bba4f8016cf207827d76bef0020e7fb.de-de.xlf") {
        $h.Delete()
    }
}
$wsDe.Range("C2").Clear()

$wsDe.Range("D2").Value = $epoch
$wsDe.Range("G2").Value = $epoch
$wsDe.Range("H2").Value = $newReason

$wsDe.Range("D3").Value = $epoch
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = $newReason

Update-Md-Hyperlink $wsDe
